# The underlying OOXML diff for this revision is a pure re-serialization
# of word/document.xml and word/styles.xml: every changed line is the same
# set of attributes on the same element, just re-ordered (alphabetically)
# by the tool that re-saved the package. No text, formatting, style
# definition, page-setup value, or any other document content actually
# changed between the two revisions.
#
# Word's COM object model only exposes document *content* (text, runs,
# paragraphs, styles, page setup, properties, ...); it has no notion of
# "XML attribute order" to replay, since Word itself re-serializes parts
# however its own writer chooses whenever it touches them. So the
# content-level edit that reproduces this revision is a no-op: open the
# document and leave every value exactly as it is.
$d = $word.ActiveDocument

# Touch the document object (establishes $d / confirms the session is
# alive) without reading or writing any property that would mark a part
# dirty, so the package is saved back out unchanged.
$null = $d.Name
